$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Talk about Bootstrap here" text box (Text Box 194 / id 15) by
# scanning shapes for the current placeholder text, rather than hard-coding
# an index, to stay robust to shape ordering.
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $t = $shape.TextFrame.TextRange.Text
        if ($t -like "*Talk about Bootstrap here*") {
            $targetShape = $shape
            break
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Build the full replacement text. Each paragraph is separated with a
# carriage return; blank lines become their own (empty) paragraph.
$para1 = "Bootstrap " + "is an open-source front-end web framework"
$para2 = ""
$para3 = "Bootstrap " + "contains CSS design templates for common HTML elements"
$para4 = ""
$para5 = "This made designing our site simple and straightforward and resulted in a clean look that is easy to view and understand"
$para6 = ""
$para7 = "Also, since " + "Bootstrap " + "is open-source, there are a variety of plugins that have been made for it"
$para8 = ""
$para9 = "These plugins allowed us to " + "add some functionality that would have taken much more " + "time and effort if we were to write the code for them on " + "our own"

$fullText = $para1 + "`r" + $para2 + "`r" + $para3 + "`r" + $para4 + "`r" + $para5 + "`r" + $para6 + "`r" + $para7 + "`r" + $para8 + "`r" + $para9

$tr.Text = $fullText

# Setting .Text above makes the whole body inherit the shape's original
# bold lead-in run formatting. Explicitly clear Bold on every span that
# must NOT be bold, leaving the "Bootstrap " lead-ins (paragraphs 1, 3, 7)
# bold as-is.
$tr.Characters(11, 41).Font.Bold = $false    # "is an open-source front-end web framework"
$tr.Characters(64, 54).Font.Bold = $false    # "contains CSS design templates for common HTML elements"
$tr.Characters(120, 120).Font.Bold = $false  # "This made designing our site..."
$tr.Characters(242, 12).Font.Bold = $false   # "Also, since "
$tr.Characters(264, 73).Font.Bold = $false   # "is open-source, there are a variety..."
$tr.Characters(339, 28).Font.Bold = $false   # "These plugins allowed us to "
$tr.Characters(367, 55).Font.Bold = $false   # "add some functionality..."
$tr.Characters(422, 57).Font.Bold = $false   # "time and effort if we were to write..."
$tr.Characters(479, 7).Font.Bold = $false    # "our own"
